$wb = $excel.ActiveWorkbook
$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# --- DBS sheet: insert a new "findGdrNum2" lookup row above the existing
# "findMainLgtseq" row (i.e. shift rows 3-5 down to 4-6, without disturbing
# the unrelated blank-cell block further down the sheet at rows 10-12). ---
for ($r = 5; $r -ge 3; $r--) {
    $destRow = $r + 1
    for ($c = 1; $c -le 3; $c++) {
        $srcCell = $wsDBS.Cells.Item($r, $c)
        $dstCell = $wsDBS.Cells.Item($destRow, $c)
        $dstCell.Value2 = $srcCell.Value2
    }
}

$wsDBS.Range("A3").Value2 = "findGdrNum2"
$wsDBS.Range("B3").Value2 = "GdrId1 = ,AND GdrId2 = ,AND GdrNum = ,AND LgtSeq ="
$wsDBS.Range("C3").Value2 = "GdrId1 ASC,GdrId2 ASC,GdrNum ASC,LgtSeq ASC"

# New row 3 picks up the same 18pt custom row height used by the header/row2.
$wsDBS.Rows.Item(3).RowHeight = 18

# --- Selection / active-sheet bookkeeping, matching the saved UI state. ---
[void]$wsDBD.Range("B11").Select()
[void]$wsDBS.Activate()
[void]$wsDBS.Range("B12").Select()
